# Scheduled market-data refresh for the Pandaemonium_Profits workbook.
# Updates the price / leve-profit columns (H:N) on the leve rows whose
# item prices moved since the last run (values sourced from the market board).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 51: A Bile Business / Shark Oil
$ws.Range("H51").Value = 3483.5
$ws.Range("I51").Value = 10001
$ws.Range("J51").Value = 2180
$ws.Range("K51").Value = 10001
$ws.Range("L51").Value = 2180
$ws.Range("M51").Value = -9517
$ws.Range("N51").Value = -3148

# Row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 3850.9644
$ws.Range("I64").Value = 3490.8948
$ws.Range("J64").Value = 4611.1113
$ws.Range("K64").Value = 3490.8948
$ws.Range("L64").Value = 4611.1113
$ws.Range("M64").Value = -3242.8948
$ws.Range("N64").Value = -5107.1113

# Row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 3850.9644
$ws.Range("I67").Value = 3490.8948
$ws.Range("J67").Value = 4611.1113
$ws.Range("K67").Value = 3490.8948
$ws.Range("L67").Value = 4611.1113
$ws.Range("M67").Value = -2632.8948
$ws.Range("N67").Value = -6327.1113

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 4201.6665
$ws.Range("I74").Value = 3860.375
$ws.Range("J74").Value = 4884.25
$ws.Range("K74").Value = 3860.375
$ws.Range("L74").Value = 4884.25
$ws.Range("M74").Value = -2924.375
$ws.Range("N74").Value = -6756.25

# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 3649.1956
$ws.Range("I76").Value = 3504.647
$ws.Range("J76").Value = 4058.75
$ws.Range("K76").Value = 3504.647
$ws.Range("L76").Value = 4058.75
$ws.Range("M76").Value = -3189.647
$ws.Range("N76").Value = -4688.75

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 4201.6665
$ws.Range("I77").Value = 3860.375
$ws.Range("J77").Value = 4884.25
$ws.Range("K77").Value = 19301.875
$ws.Range("L77").Value = 24421.25
$ws.Range("M77").Value = -14621.875
$ws.Range("N77").Value = -33781.25

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 3649.1956
$ws.Range("I79").Value = 3504.647
$ws.Range("J79").Value = 4058.75
$ws.Range("K79").Value = 3504.647
$ws.Range("L79").Value = 4058.75
$ws.Range("M79").Value = -2412.647
$ws.Range("N79").Value = -6242.75

# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 9562.637000000001
$ws.Range("I125").Value = 221.85715
$ws.Range("J125").Value = 25909
$ws.Range("K125").Value = 1996.71435
$ws.Range("L125").Value = 233181
$ws.Range("M125").Value = 463.28565
$ws.Range("N125").Value = -238101


$ws = $wb.Worksheets.Item("ARM")

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 10000.857
$ws.Range("I63").Value = 2250
$ws.Range("J63").Value = 13101.2
$ws.Range("K63").Value = 2250
$ws.Range("L63").Value = 13101.2
$ws.Range("M63").Value = -1564
$ws.Range("N63").Value = -14473.2

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 10000.857
$ws.Range("I66").Value = 2250
$ws.Range("J66").Value = 13101.2
$ws.Range("K66").Value = 11250
$ws.Range("L66").Value = 65506
$ws.Range("M66").Value = -7818
$ws.Range("N66").Value = -72370


$ws = $wb.Worksheets.Item("BSM")

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 4902.8335
$ws.Range("I105").Value = 4739.636
$ws.Range("J105").Value = 5501.222
$ws.Range("K105").Value = 4739.636
$ws.Range("L105").Value = 5501.222
$ws.Range("M105").Value = -2992.636
$ws.Range("N105").Value = -8995.222


$ws = $wb.Worksheets.Item("CRP")

# Row 6: Got Your Back / Square Maple Shield
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 5000
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -5226

# Row 7: Gridania's Got Talent / Maple Lumber
$ws.Range("H7").Value = 88.888885
$ws.Range("I7").Value = 87.5
$ws.Range("J7").Value = 90
$ws.Range("K7").Value = 87.5
$ws.Range("L7").Value = 90
$ws.Range("M7").Value = 25.5
$ws.Range("N7").Value = -316

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3630.2666
$ws.Range("I31").Value = 1440
$ws.Range("J31").Value = 4898.316
$ws.Range("K31").Value = 1440
$ws.Range("L31").Value = 4898.316
$ws.Range("M31").Value = -1145
$ws.Range("N31").Value = -5488.316

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3630.2666
$ws.Range("I34").Value = 1440
$ws.Range("J34").Value = 4898.316
$ws.Range("K34").Value = 1440
$ws.Range("L34").Value = 4898.316
$ws.Range("M34").Value = -1238
$ws.Range("N34").Value = -5302.316

# Row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 1376.8462
$ws.Range("I107").Value = 1444.1818
$ws.Range("J107").Value = 1006.5
$ws.Range("K107").Value = 1444.1818
$ws.Range("L107").Value = 1006.5
$ws.Range("M107").Value = 475.8181999999999
$ws.Range("N107").Value = -4846.5


$ws = $wb.Worksheets.Item("CUL")

# Row 55: Pagan Pastries / Pastry Fish
$ws.Range("H55").Value = 8143.2856
$ws.Range("I55").Value = 5004
$ws.Range("K55").Value = 15012
$ws.Range("M55").Value = -14835

# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 1278.3
$ws.Range("I68").Value = 797.5
$ws.Range("J68").Value = 1398.5
$ws.Range("K68").Value = 2392.5
$ws.Range("L68").Value = 4195.5
$ws.Range("M68").Value = -1581.5
$ws.Range("N68").Value = -5817.5

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 1278.3
$ws.Range("I71").Value = 797.5
$ws.Range("J71").Value = 1398.5
$ws.Range("K71").Value = 7177.5
$ws.Range("L71").Value = 12586.5
$ws.Range("M71").Value = -3121.5
$ws.Range("N71").Value = -20698.5

# Row 111: Soup for the Soldier / Broad Bean Soup
$ws.Range("H111").Value = 2870.75
$ws.Range("I111").Value = 1151
$ws.Range("J111").Value = 8030
$ws.Range("K111").Value = 3453
$ws.Range("L111").Value = 24090
$ws.Range("M111").Value = -386
$ws.Range("N111").Value = -30224

# Row 112: Sweet Tooth / Caramels
$ws.Range("H112").Value = 4357.143
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 5000
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 15000
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -17216

# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 775.46155
$ws.Range("J113").Value = 759.6667
$ws.Range("L113").Value = 2279.0001
$ws.Range("N113").Value = -6619.0001


$ws = $wb.Worksheets.Item("GSM")

# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 5254.3687
$ws.Range("J70").Value = 5846.9585
$ws.Range("L70").Value = 5846.9585
$ws.Range("N70").Value = -6386.9585

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 5254.3687
$ws.Range("J73").Value = 5846.9585
$ws.Range("L73").Value = 5846.9585
$ws.Range("N73").Value = -7718.9585


$ws = $wb.Worksheets.Item("LTW")

# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 3857.8696
$ws.Range("I40").Value = 3670.2778
$ws.Range("J40").Value = 4533.2
$ws.Range("K40").Value = 3670.2778
$ws.Range("L40").Value = 4533.2
$ws.Range("M40").Value = -3534.2778
$ws.Range("N40").Value = -4805.2

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 628.5714
